# Update "paises.xlsx" countries & provincias Spain data.
# - Re-labels the 4 country-pairs whose shared-string ordering swapped
#   (Egipto/Suecia, Tanzania/Surinam, Laos/Santa Lucia, Fiyi/Dominica)
# - Refreshes the numeric case counters for the affected rows
# - Bumps the "Datos actualizados" timestamp banner

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Country label swaps (rows identified by their position in the
#        "Casos totales" ranking table) -------------------------------
$ws.Range("A26").Value = "Suecia"
$ws.Range("A27").Value = "Egipto"

$ws.Range("A155").Value = "Surinam"
$ws.Range("A156").Value = "Tanzania"

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# --- 2) Updated numeric data (Casos totales, Nuevos casos, Casos
#        activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2716543
$ws.Range("C4").Value = 34732
$ws.Range("D4").Value = 1127575
$ws.Range("E4").Value = 1459005
$ws.Range("G4").Value = 564
$ws.Range("H4").Value = 129963

# Row 5 - Brasil
$ws.Range("B5").Value = 1402041
$ws.Range("C5").Value = 31553
$ws.Range("E5").Value = 584985
$ws.Range("G5").Value = 1209
$ws.Range("H5").Value = 59594

# Row 10 - Peru
$ws.Range("B10").Value = 285213
$ws.Range("C10").Value = 2848
$ws.Range("D10").Value = 174535
$ws.Range("E10").Value = 101001
$ws.Range("G10").Value = 173
$ws.Range("H10").Value = 9677

# Row 17 - Alemania
$ws.Range("B17").Value = 195797
$ws.Range("C17").Value = 405
$ws.Range("E17").Value = 7645
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 9052

# Row 26 - now Suecia
$ws.Range("B26").Value = 68451
$ws.Range("C26").Value = 198
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 5333

# Row 27 - now Egipto
$ws.Range("B27").Value = 68311
$ws.Range("C27").Value = 1557
$ws.Range("D27").Value = 18460
$ws.Range("E27").Value = 46898
$ws.Range("G27").Value = 81
$ws.Range("H27").Value = 2953

# Row 46 - Suiza
$ws.Range("D46").Value = 29200
$ws.Range("E46").Value = 551

# Row 49 - Barein
$ws.Range("B49").Value = 26758
$ws.Range("C49").Value = 519
$ws.Range("D49").Value = 21331
$ws.Range("E49").Value = 5340

# Row 74 - Uzbekistan
$ws.Range("B74").Value = 8503
$ws.Range("C74").Value = 281
$ws.Range("E74").Value = 2795

# Row 89 - Bulgaria
$ws.Range("B89").Value = 4989
$ws.Range("C89").Value = 158
$ws.Range("D89").Value = 2676
$ws.Range("E89").Value = 2083
$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 230

# Row 134 - Ruanda
$ws.Range("B134").Value = 1025
$ws.Range("C134").Value = 24
$ws.Range("D134").Value = 447
$ws.Range("E134").Value = 576

# Row 151 - Togo
$ws.Range("B151").Value = 650
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 402
$ws.Range("E151").Value = 234

# Row 155 - now Surinam
$ws.Range("B155").Value = 515
$ws.Range("C155").Value = 14
$ws.Range("D155").Value = 227
$ws.Range("E155").Value = 275
$ws.Range("H155").Value = 13

# Row 156 - now Tanzania
$ws.Range("B156").Value = 509
$ws.Range("D156").Value = 183
$ws.Range("E156").Value = 305
$ws.Range("H156").Value = 21

# --- 3) Timestamp banner -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 23:48"
